# Add the new weekly ranking sheet "magapoke_2025-12-31" at the end of the
# workbook, matching the layout/style of the existing magapoke_* sheets.

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
$firstSheet = $wb.Worksheets.Item(1)
$lastSheet = $wb.Worksheets.Item($sheetCount)

$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "magapoke_2025-12-31"

# Match the page margins used by every other magapoke_* sheet in the
# workbook (0.75in/0.75in/1in/1in/0.5in/0.5in -> PageSetup uses points).
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Header row (values + reuse the existing header formatting: bold, centered,
# thin border, matching every other magapoke_* sheet in the workbook).
$ws.Range("A1").Value = "rank"
$ws.Range("B1").Value = "title"

$firstSheet.Range("A1:B1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

# Ranking rows (rank 1..44, title)
$titles = @(
    '黒月のイェルクナハト',
    'スルガメテオ',
    'ドリーム☆ジャンボ☆ガール',
    '黄昏町プリズナーズ',
    'アイドラトリィ',
    'K-9~警視庁公安部公安第9課異能対策係~',
    'せいぶつ部の田辺くん',
    'ハードワーカー中田',
    'ゼロとヒャク',
    '篝家の８兄弟',
    '生きたがりの人狼',
    'ともだちづくり',
    'ナキナギ',
    '追放されなかった男　～二度目の人生は土下座から始まりました～',
    'ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜',
    '普通の本はありません！',
    'お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！',
    '春くらり',
    '夜鐘のキト',
    '皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～',
    'MYS',
    '歪みの虜',
    '限界集落を脱村した錬金術士、都会で"最強"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～',
    'その青春',
    '卒業アルバムの彼女たち',
    '異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～',
    '屋根の下のアルテミス',
    '君が監督！',
    'きゃわるり方程式',
    '明智ナンバーワン',
    'じゅーくぼっくす',
    'ハプスブルク家の華麗なる受難',
    '鳴るさんだぁ',
    '平成転生',
    'GURU',
    'JK Biker',
    '永久のユウグレ',
    'ナマイキ旭ちゃんをわからせたい',
    '白銀のキュイジーヌ～明治外交官の料理人～',
    '眠れる森のレガ',
    '〈小市民〉 春期限定いちごタルト事件',
    '花子狩り',
    '人生逆転ダンジョン',
    'イエティ、とある日々'
)

for ($i = 0; $i -lt $titles.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $titles[$i]
}
